$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2, shifting existing row 2 (and below) down to row 3.
# Rows.Insert() copies formatting from the row above (the bold header row) onto
# the new row, so strip that back off to match the unformatted data rows.
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).ClearFormats()

# Fill in the new job posting in row 2
$ws.Range("A2").Value = "Oracle Careers"
$ws.Range("B2").Value = "Agri Finance and R&D Specialist (Open to external applicants)"
$ws.Range("C2").Value = "New Delhi, India"
# Leading apostrophe forces this to stay plain text instead of being
# auto-converted to a date serial number.
$ws.Range("D2").Value = "'01/21/2026"
$ws.Range("E2").Formula = '=HYPERLINK("https://estm.fa.em2.oraclecloud.com/hcmUI/CandidateExperience/en/sites/CX_1/job/31390/?location=India&locationId=300000000440677&locationLevel=country&mode=location", "Apply")'

# Setting the D2 text triggers an auto-applied "Text" number format (style),
# which the target row doesn't carry - clear it again so row 2 ends up with
# the same (unstyled) look as the other data rows.
$ws.Rows.Item(2).ClearFormats()
